$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 55
$ws.Range("H55").Value = 315
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 315
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 315
$ws.Range("N55").Value = -743
$ws.Range("M55").ClearContents()

# ALC row 87
$ws.Range("H87").Value = 41140
$ws.Range("J87").Value = 41140
$ws.Range("L87").Value = 41140
$ws.Range("N87").Value = -43636

# ALC row 90
$ws.Range("H90").Value = 41140
$ws.Range("J90").Value = 41140
$ws.Range("L90").Value = 123420
$ws.Range("N90").Value = -135900

# ALC row 111
$ws.Range("H111").Value = 1048.75
$ws.Range("I111").Value = 957
$ws.Range("K111").Value = 2871
$ws.Range("M111").Value = 196

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 2816.1667
$ws.Range("I32").Value = 1549.375
$ws.Range("K32").Value = 1549.375
$ws.Range("M32").Value = -1262.375

# ARM row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# ARM row 122
$ws.Range("H122").Value = 7178
$ws.Range("I122").Value = 8232.375
$ws.Range("J122").Value = 5491
$ws.Range("K122").Value = 24697.125
$ws.Range("L122").Value = 16473
$ws.Range("M122").Value = -22247.125
$ws.Range("N122").Value = -21373

# ARM row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# ARM row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ARM row 140
$ws.Range("H140").Value = 146666.33
$ws.Range("J140").Value = 146666.33
$ws.Range("L140").Value = 146666.33
$ws.Range("N140").Value = -157026.33

# ARM row 141
$ws.Range("H141").Value = 195000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 195000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 195000
$ws.Range("N141").Value = -205360
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20
$ws.Range("H20").Value = 1126.2222
$ws.Range("I20").Value = 1079.5
$ws.Range("J20").Value = 1500
$ws.Range("K20").Value = 1079.5
$ws.Range("L20").Value = 1500
$ws.Range("M20").Value = -832.5
$ws.Range("N20").Value = -1994

# BSM row 64
$ws.Range("H64").Value = 860.8333
$ws.Range("I64").Value = 895
$ws.Range("J64").Value = 843.75
$ws.Range("K64").Value = 895
$ws.Range("L64").Value = 843.75
$ws.Range("M64").Value = -670
$ws.Range("N64").Value = -1293.75

# BSM row 67
$ws.Range("H67").Value = 860.8333
$ws.Range("I67").Value = 895
$ws.Range("J67").Value = 843.75
$ws.Range("K67").Value = 895
$ws.Range("L67").Value = 843.75
$ws.Range("M67").Value = -115
$ws.Range("N67").Value = -2403.75

# BSM row 75
$ws.Range("H75").Value = 4300
$ws.Range("I75").Value = 4300
$ws.Range("K75").Value = 4300
$ws.Range("M75").Value = -3364

# BSM row 78
$ws.Range("H78").Value = 4300
$ws.Range("I78").Value = 4300
$ws.Range("K78").Value = 12900
$ws.Range("M78").Value = -8220

# BSM row 105
$ws.Range("H105").Value = 20228.889
$ws.Range("I105").Value = 20228.889
$ws.Range("K105").Value = 20228.889
$ws.Range("M105").Value = -18481.889

# BSM row 134
$ws.Range("H134").Value = 1806.5333
$ws.Range("I134").Value = 1188.2727
$ws.Range("K134").Value = 3564.8181
$ws.Range("M134").Value = -1029.8181

# BSM row 140
$ws.Range("H140").Value = 107500
$ws.Range("J140").Value = 107500
$ws.Range("L140").Value = 107500
$ws.Range("N140").Value = -117860

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16
$ws.Range("H16").Value = 586.9167
$ws.Range("I16").Value = 445.22223
$ws.Range("J16").Value = 1012
$ws.Range("K16").Value = 445.22223
$ws.Range("L16").Value = 1012
$ws.Range("M16").Value = -158.22223
$ws.Range("N16").Value = -1586

# CRP row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

# CRP row 113
$ws.Range("H113").Value = 586.9167
$ws.Range("I113").Value = 445.22223
$ws.Range("J113").Value = 1012
$ws.Range("K113").Value = 445.22223
$ws.Range("L113").Value = 1012
$ws.Range("M113").Value = 1724.77777
$ws.Range("N113").Value = -5352

$ws = $wb.Worksheets.Item("CUL")
# CUL row 11
$ws.Range("H11").Value = 100
$ws.Range("I11").Value = 100
$ws.Range("K11").Value = 300
$ws.Range("M11").Value = -160

# CUL row 12
$ws.Range("H12").Value = 44
$ws.Range("J12").Value = 54.75
$ws.Range("L12").Value = 164.25
$ws.Range("N12").Value = -510.25

# CUL row 55
$ws.Range("H55").Value = 2539.2856
$ws.Range("I55").Value = 2633.3333
$ws.Range("J55").Value = 2468.75
$ws.Range("K55").Value = 7899.999899999999
$ws.Range("L55").Value = 7406.25
$ws.Range("M55").Value = -7722.999899999999
$ws.Range("N55").Value = -7760.25

# CUL row 139
$ws.Range("H139").Value = 3912.3635
$ws.Range("I139").Value = 3912.3635
$ws.Range("K139").Value = 11737.0905
$ws.Range("M139").Value = -6597.0905

# CUL row 140
$ws.Range("H140").Value = 1347.1428
$ws.Range("I140").Value = 1347.1428
$ws.Range("K140").Value = 4041.4284
$ws.Range("M140").Value = 1138.5716

$ws = $wb.Worksheets.Item("GSM")
# GSM row 141
$ws.Range("H141").Value = 57500
$ws.Range("J141").Value = 57500
$ws.Range("L141").Value = 57500
$ws.Range("N141").Value = -67860

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705

# LTW row 27
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893

# LTW row 55
$ws.Range("H55").Value = 3108
$ws.Range("I55").Value = 615.2
$ws.Range("J55").Value = 7262.6665
$ws.Range("K55").Value = 615.2
$ws.Range("L55").Value = 7262.6665
$ws.Range("M55").Value = -442.2
$ws.Range("N55").Value = -7608.6665

# LTW row 82
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# LTW row 85
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# LTW row 99
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 40000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 40000
$ws.Range("N99").Value = -45990
$ws.Range("M99").ClearContents()

# LTW row 102
$ws.Range("H102").Value = 70200
$ws.Range("J102").Value = 70200
$ws.Range("L102").Value = 70200
$ws.Range("N102").Value = -76690

# LTW row 138
$ws.Range("H138").Value = 70000
$ws.Range("I138").Value = 70000
$ws.Range("K138").Value = 70000
$ws.Range("M138").Value = -64860

$ws = $wb.Worksheets.Item("WVR")
# WVR row 102
$ws.Range("H102").Value = 89999
$ws.Range("J102").Value = 89999
$ws.Range("L102").Value = 89999
$ws.Range("N102").Value = -96489

# WVR row 132
$ws.Range("H132").Value = 3332.5
$ws.Range("I132").Value = 2998.75
$ws.Range("K132").Value = 8996.25
$ws.Range("M132").Value = -6466.25

# WVR row 135
$ws.Range("H135").Value = 51250
$ws.Range("J135").Value = 51250
$ws.Range("L135").Value = 51250
$ws.Range("N135").Value = -61390

# WVR row 137
$ws.Range("H137").Value = 40000
$ws.Range("J137").Value = 40000
$ws.Range("L137").Value = 40000
$ws.Range("N137").Value = -50200

# WVR row 139
$ws.Range("H139").Value = 60000
$ws.Range("J139").Value = 60000
$ws.Range("L139").Value = 60000
$ws.Range("N139").Value = -70280
